$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (102 and 103) to the feed logs sheet,
# matching the existing data layout: run_id, rss_url_id, date, response, item_count

$ws.Cells.Item(102, 1).Value = 101
$ws.Cells.Item(102, 2).Value = 1
$ws.Cells.Item(102, 3).Value = "2024-06-17 02:05:25"
$ws.Cells.Item(102, 4).Value = 200
$ws.Cells.Item(102, 5).Value = 4

$ws.Cells.Item(103, 1).Value = 102
$ws.Cells.Item(103, 2).Value = 2
$ws.Cells.Item(103, 3).Value = "2024-06-17 02:05:26"
$ws.Cells.Item(103, 4).Value = 200
$ws.Cells.Item(103, 5).Value = 0
